$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")

# --- Settings sheet: replace the old QueueName/ProcessQueue row and add the
# new Orchestrator-related configuration rows (2-5). Values are written in
# the same order they were typed in the original editing session so that
# the shared-string table is rebuilt with a matching append order.
$settings.Range("B4").Value = "fantastic"
$settings.Range("B2").Value = "https://demo.uipath.com"
$settings.Range("B5").Value = "KibanaDemoQueue"
$settings.Range("B3").Value = "demo.uipath.com_credentials"
$settings.Range("A2").Value = "OrchestratorURL"
$settings.Range("C5").Value = "Orchestrator Queue Name. Be sure to match this name with the one from the server."
$settings.Range("A3").Value = "OrchestratorCredentialName"
$settings.Range("A4").Value = "OrchestratorTenancyName"
$settings.Range("A5").Value = "OrchestratorQueueName"
$settings.Range("C2").Value = "The URL of your orchestrator server. This property is used only if you are using a Queue to store your Transaction Items."
$settings.Range("C3").Value = "The name of Orchestrator credentials. This should be stored in Windows Credential manager. This property is used only if you are using a Queue to store your Transaction Items."
$settings.Range("C4").Value = "The name of the Orchestrator tenant.  This property is used only if you are using a Queue to store your Transaction Items."

# --- Constants sheet: re-enter the same log-message text for rows 21-23
# (content unchanged, but re-typing lines these cells back up against the
# already-existing shared strings instead of the removed ones).
$constants.Range("B21").Value = "Transaction Successful."
$constants.Range("B22").Value = "Business rule exception."
$constants.Range("B23").Value = "System exception."

# --- Selection / active sheet bookkeeping to match the final view state.
# Constants keeps a stale selection at A29 (below its data), while Settings
# becomes the active/visible tab with A2 selected - selecting Settings last
# is what makes it the active sheet.
$constants.Range("A29").Select() | Out-Null
$settings.Range("A2").Select() | Out-Null
